$wb = $excel.ActiveWorkbook

$wsTestCases = $wb.Worksheets.Item("TestCases")

# Add the new test case row (row 21) on the TestCases sheet.
$wsTestCases.Cells.Item(21, 1).Value = "GK_016_Test"
$wsTestCases.Cells.Item(21, 2).Value = "Verify user is able to add all vegitables."

# Make TestCases the active sheet / tab, with the view scrolled down and
# the given cell selected, matching the author's final on-screen state.
$wsTestCases.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 7
$wsTestCases.Range("C20").Select() | Out-Null
